# Update the "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages rebuild at 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> { row number -> new F-column value }
$updates = @{
    "展览" = @{
        3  = 232
        4  = 18
        5  = 6751
        7  = 5
        8  = 435
        9  = 144
        10 = 6300
        13 = 1271
        15 = 102
        19 = 370
        22 = 4653
        23 = 65
        24 = 47
        25 = 109
        27 = 82
    }
    "全部类型" = @{
        3  = 232
        4  = 18
        5  = 6751
        7  = 5
        8  = 435
        9  = 144
        10 = 6300
        13 = 1271
        15 = 102
        19 = 370
        22 = 4653
        24 = 65
        25 = 47
        26 = 109
        28 = 82
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($rowNum in $rows.Keys) {
        $ws.Cells.Item($rowNum, 6).Value = $rows[$rowNum]
    }
}
